# Apply "Añadidos precios del día 26 al script inicial" changes
# to the SCRAP_PRECIOS sheet (prices scraping log).

$wb = $excel.ActiveWorkbook
$wsScrap = $wb.Worksheets.Item("SCRAP_PRECIOS")

# --- Update day-25 price row (row 3) to its day-26 (26 Jan 2021) value ---
$wsScrap.Range("D3").Value = 144.47999999999999
$wsScrap.Range("E3").Value = 44222

# --- Insert a new row for the 26 Jan 2021 price of IE00B4WXT857 ---
$wsScrap.Rows("25:25").Insert()

$wsScrap.Range("C25").Value = "IE00B4WXT857"
$wsScrap.Range("D25").Value = 13.75
$wsScrap.Range("E25").Value = 44222

# --- Correct the price that is now on row 27 (was row 26) ---
$wsScrap.Range("D27").Value = 13.27

# --- Make SCRAP_PRECIOS the active sheet/tab ---
$wsScrap.Activate()
